$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "Cost Summary"
$wb.Worksheets.Item(8).Name = "Operating and Support Cost"

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("C23").Select()

$ws8 = $wb.Worksheets.Item(8)
$ws8.Activate()
$ws8.Range("E16").Select()
